$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 296, pushing existing
# rows 296-299 down to 298-301 (their original values are preserved
# automatically by the insert operation).
$ws.Rows.Item(296).Insert()
$ws.Rows.Item(296).Insert()

# New row 296: weekly update for "Provincia de Linares"
$ws.Cells.Item(296,1).Value = 6
$ws.Cells.Item(296,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(296,3).Value = "Metropolitana"
$ws.Cells.Item(296,4).Value = 44595
$ws.Cells.Item(296,5).Value = 13
$ws.Cells.Item(296,6).Value = "Fruta"
$ws.Cells.Item(296,7).Value = 100101
$ws.Cells.Item(296,8).Value = "Berries"
$ws.Cells.Item(296,9).Value = 100101001
$ws.Cells.Item(296,10).Value = "Arándano (blue)"
$ws.Cells.Item(296,11).Value = "Sin especificar"
$ws.Cells.Item(296,12).Value = "Primera"
$ws.Cells.Item(296,13).Value = 200
$ws.Cells.Item(296,14).Value = 4000
$ws.Cells.Item(296,15).Value = 4000
$ws.Cells.Item(296,16).Value = 4000
$ws.Cells.Item(296,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(296,18).Value = "Provincia de Linares"
$ws.Cells.Item(296,19).Value = 2000
$ws.Cells.Item(296,20).Value = 2

# New row 297: weekly update for "Región de O'Higgins"
$ws.Cells.Item(297,1).Value = 6
$ws.Cells.Item(297,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(297,3).Value = "Metropolitana"
$ws.Cells.Item(297,4).Value = 44595
$ws.Cells.Item(297,5).Value = 13
$ws.Cells.Item(297,6).Value = "Fruta"
$ws.Cells.Item(297,7).Value = 100101
$ws.Cells.Item(297,8).Value = "Berries"
$ws.Cells.Item(297,9).Value = 100101001
$ws.Cells.Item(297,10).Value = "Arándano (blue)"
$ws.Cells.Item(297,11).Value = "Sin especificar"
$ws.Cells.Item(297,12).Value = "Primera"
$ws.Cells.Item(297,13).Value = 340
$ws.Cells.Item(297,14).Value = 4000
$ws.Cells.Item(297,15).Value = 4500
$ws.Cells.Item(297,16).Value = 4250
$ws.Cells.Item(297,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(297,18).Value = "Región de O'Higgins"
$ws.Cells.Item(297,19).Value = 2125
$ws.Cells.Item(297,20).Value = 2
